# Updates a handful of "DD_#"/"SD_#" label cells on the ten_lists sheet,
# replacing old list labels (A1, A3, A7, A8, A11) with the new set of
# inflatable lists (B1, B3, B7, B8, B11) added for the Older Adults project.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

$ws.Range("H5").Value  = "B1"
$ws.Range("H12").Value = "B3"
$ws.Range("J15").Value = "B11"
$ws.Range("D19").Value = "B3"
$ws.Range("J26").Value = "B3"
$ws.Range("B28").Value = "B7"
$ws.Range("H28").Value = "B8"

$ws.Range("J27").Select()
